$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update table contents (ItemName/Ingredients/Allergens/LocalIngredients/Diet/nutritionLabel/LeaveEmpty) ---
$ws.Range("A1").Value = "ItemName"
$ws.Range("B1").Value = "Ingredients"
$ws.Range("C1").Value = "Allergens"
$ws.Range("D1").Value = "LocalIngredients"
$ws.Range("E1").Value = "Diet"
$ws.Range("F1").Value = "nutritionLabel"
$ws.Range("G1").Value = "LeaveEmpty"

$ws.Range("A2").Value = "Bacon Breakfast Burrito"
$ws.Range("B2").Value = "Eggs / Bacon / Tater Tots /  Salsa / Nacho Cheese"
$ws.Range("C2").Value = "Wheat, gluten, milk, sulphites, egg."
$ws.Range("D2").Value = "Castle Cheese Nacho Blend"
$ws.Range("E2").Value = "GFO"
$ws.Range("F2").Value = "Bacon_Breakfast_Burrito"

$ws.Range("A3").Value = "Sausage Breakfast Burrito"
$ws.Range("B3").Value = "Eggs / Pork Sausage / Tater Tots /  Salsa / Nacho Cheese"
$ws.Range("C3").Value = "Wheat, gluten, milk, sulphites, egg."
$ws.Range("D3").Value = " Castle Cheese Nacho Blend"
$ws.Range("E3").Value = "GFO"
$ws.Range("F3").Value = "Bacon_Breakfast_Burrito"

$ws.Range("A4").Value = "Veggie Breakfast Burrito"
$ws.Range("B4").Value = "Eggs / Spinach / Tater Tots /  Salsa / Nacho Cheese"
$ws.Range("C4").Value = "Wheat, gluten, milk, sulphites, egg."
$ws.Range("D4").Value = " Castle Cheese Nacho Blend"
$ws.Range("E4").Value = "VGN, GFO"
$ws.Range("F4").Value = "Falafel_Wrap"

$ws.Range("A5").Value = "Breakfast Special"
$ws.Range("B5").Value = "Eggs / Bacon or Sausage / Tater Tots / Croissant"
$ws.Range("C5").Value = "Wheat, gluten, milk, sulphites, egg."
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "Bacon_Breakfast_Burrito"

$ws.Range("A6").Value = "Vegetarian Breakfast Special"
$ws.Range("B6").Value = "Eggs / Vegan Sausage / Tater Tots / Croissant"
$ws.Range("C6").Value = "Wheat, gluten, milk, sulphites, egg."
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "VEG"
$ws.Range("F6").Value = "Falafel_Wrap"

$ws.Range("A7").Value = "Lumberjack Bowl"
$ws.Range("B7").Value = "Tater Tots / Sausage  and  Bacon / Cheddar Cheese / Diced Tomato / Green Onion / Hollandaise"
$ws.Range("C7").Value = "Wheat, gluten, milk, sulphites, egg."
$ws.Range("D7").Value = "Fresh Start Tomato, Castle Cheese Nacho Blend"
$ws.Range("E7").Value = "BC"
$ws.Range("F7").Value = "Bacon_Breakfast_Burrito"

$ws.Range("A8").Value = "Vegetarian Lumberjack Bowl"
$ws.Range("B8").Value = "Tater Tots / Vegan Sausage / Cheddar Cheese / Diced Tomato / Green Onion / Hollandaise"
$ws.Range("C8").Value = "Wheat, gluten, milk, sulphites, egg."
$ws.Range("D8").Value = "Fresh Start Tomato, Castle Cheese Nacho Blend"
$ws.Range("E8").Value = "VEG, BC"
$ws.Range("F8").Value = "Falafel_Wrap"

$ws.Range("A9").Value = "Breakfast Waffles"
$ws.Range("B9").Value = "Waffles / Berry Compote / Whipped Cream"
$ws.Range("C9").Value = "Wheat, milk, egg."
$ws.Range("D9").Value = "Patisserie LeBeau Waffle"
$ws.Range("E9").Value = "VEG"
$ws.Range("F9").Value = "Apple_Turnover"


# --- Remove the wrap-text style previously applied to B3 (back to Normal/default style) ---
$ws.Range("B3").Style = "Normal"

# --- Update the active selection/cell shown when the workbook is opened ---
$ws.Range("F9").Select()
